# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting a refreshed scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 887
$ws1.Range("F4").Value = 790
$ws1.Range("F6").Value = 448
$ws1.Range("F9").Value = 1287
$ws1.Range("F10").Value = 716
$ws1.Range("F11").Value = 415
$ws1.Range("F12").Value = 546
$ws1.Range("F14").Value = 38
$ws1.Range("F15").Value = 981
$ws1.Range("F17").Value = 407
$ws1.Range("F19").Value = 92
$ws1.Range("F20").Value = 587
$ws1.Range("F22").Value = 638
$ws1.Range("F24").Value = 1011
$ws1.Range("F25").Value = 15

# --- Sheet "演出" (Show) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 243

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 887
$ws4.Range("F6").Value = 790
$ws4.Range("F8").Value = 448
$ws4.Range("F11").Value = 1287
$ws4.Range("F12").Value = 716
$ws4.Range("F15").Value = 415
$ws4.Range("F16").Value = 546
$ws4.Range("F19").Value = 38
$ws4.Range("F20").Value = 981
$ws4.Range("F23").Value = 407
$ws4.Range("F25").Value = 92
$ws4.Range("F26").Value = 243
$ws4.Range("F28").Value = 587
$ws4.Range("F34").Value = 638
$ws4.Range("F36").Value = 1012
$ws4.Range("F37").Value = 15
